$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new "season record" columns.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, border, centered) used by the rest of row 1.
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

# Fill in the season record for every player row (2 through 42) with the
# team's actual win/loss/tie totals for the season.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 29).Value = 74
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 0
}
